$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8958974811491487
$ws.Range("B3").Value = 0.855131877062405
$ws.Range("B4").Value = 0.8305502884286398
$ws.Range("B5").Value = 0.8206462812852351
$ws.Range("B6").Value = 0.8190085758407406
$ws.Range("B7").Value = 0.8304162607301748
$ws.Range("B8").Value = 0.8817485895493462
$ws.Range("B9").Value = 0.9859608992144047
$ws.Range("B10").Value = 1.064685264047455
$ws.Range("B11").Value = 1.100967605305698
$ws.Range("B12").Value = 1.114774171865974
$ws.Range("B13").Value = 1.111797696747885
$ws.Range("B14").Value = 1.10210213464444
$ws.Range("B15").Value = 1.096172066933036
$ws.Range("B16").Value = 1.062323556841193
$ws.Range("B17").Value = 1.041678753935969
$ws.Range("B18").Value = 1.029848692950253
$ws.Range("B19").Value = 1.025850855157017
$ws.Range("B20").Value = 1.04387184807328
$ws.Range("B21").Value = 1.104948134699782
$ws.Range("B22").Value = 1.145256715723576
$ws.Range("B23").Value = 1.12370755478031
$ws.Range("B24").Value = 1.042880228980266
$ws.Range("B25").Value = 0.9573891481026067
$ws.Range("C2").Value = 0.09212820554067491
$ws.Range("C3").Value = 0.09088503284039007
$ws.Range("C4").Value = 0.09010933337058447
$ws.Range("C5").Value = 0.08979011341237708
$ws.Range("C6").Value = 0.0897369189344559
$ws.Range("C7").Value = 0.09010504087667925
$ws.Range("C8").Value = 0.09170212918362353
$ws.Range("C9").Value = 0.09473592698562072
$ws.Range("C10").Value = 0.09690551159251015
$ws.Range("C11").Value = 0.09787973323489041
$ws.Range("C12").Value = 0.09824681907650756
$ws.Range("C13").Value = 0.09816784199144024
$ws.Range("C14").Value = 0.09790997029728743
$ws.Range("C15").Value = 0.09775177803206248
$ws.Range("C16").Value = 0.09684158771406715
$ws.Range("C17").Value = 0.09627995431488046
$ws.Range("C18").Value = 0.09595571910285372
$ws.Range("C19").Value = 0.09584573277419395
$ws.Range("C20").Value = 0.09633986519062177
$ws.Range("C21").Value = 0.09798576314955199
$ws.Range("C22").Value = 0.09905078176091564
$ws.Range("C23").Value = 0.09848333755931549
$ws.Range("C24").Value = 0.09631278367395879
$ws.Range("C25").Value = 0.09392565892922988
$ws.Range("E2").Value = 0.09849383471637552
$ws.Range("E3").Value = 0.09779839334300888
$ws.Range("E4").Value = 0.09741839317298684
$ws.Range("E5").Value = 0.09727537532556596
$ws.Range("E6").Value = 0.09725234261381743
$ws.Range("E7").Value = 0.0974164164421758
$ws.Range("E8").Value = 0.0982442999606512
$ws.Range("E9").Value = 0.1002402346856961
$ws.Range("E10").Value = 0.1019333408760588
$ws.Range("E11").Value = 0.1027527457792203
$ws.Range("E12").Value = 0.1030700975735215
$ws.Range("E13").Value = 0.1030014363765837
$ws.Range("E14").Value = 0.1027787130723006
$ws.Range("E15").Value = 0.1026432077367296
$ws.Range("E16").Value = 0.1018807796926744
$ws.Range("E17").Value = 0.1014256465041896
$ws.Range("E18").Value = 0.1011684979479099
$ws.Range("E19").Value = 0.1010822278549739
$ws.Range("E20").Value = 0.1014736169167811
$ws.Range("E21").Value = 0.1028439407875368
$ws.Range("E22").Value = 0.1037806766186193
$ws.Range("E23").Value = 0.1032769625641627
$ws.Range("E24").Value = 0.1014519154332838
$ws.Range("E25").Value = 0.09966044550620268
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G2").Value = 1.276085621342546
$ws.Range("G3").Value = 1.278688799162438
$ws.Range("G4").Value = 1.281007804894259
$ws.Range("G5").Value = 1.282133840922782
$ws.Range("G6").Value = 1.282331746148543
$ws.Range("G7").Value = 1.281022258310571
$ws.Range("G8").Value = 1.276833510328899
$ws.Range("G9").Value = 1.274347301887445
$ws.Range("G10").Value = 1.276028310658361
$ws.Range("G11").Value = 1.277558239926819
$ws.Range("G12").Value = 1.278247880633415
$ws.Range("G13").Value = 1.27809444499924
$ws.Range("G14").Value = 1.277612765273204
$ws.Range("G15").Value = 1.277332092965153
$ws.Range("G16").Value = 1.275943745953995
$ws.Range("G17").Value = 1.275288211648487
$ws.Range("G18").Value = 1.274983177645623
$ws.Range("G19").Value = 1.274892258894113
$ws.Range("G20").Value = 1.27535053960905
$ws.Range("G21").Value = 1.277751251050617
$ws.Range("G22").Value = 1.279963224735326
$ws.Range("G23").Value = 1.278723737827605
$ws.Range("G24").Value = 1.275322137370992
$ws.Range("G25").Value = 1.274405053638674
$ws.Range("H2").Value = 1.220571132696733
$ws.Range("H3").Value = 1.226812232622564
$ws.Range("H4").Value = 1.231151867882531
$ws.Range("H5").Value = 1.233047959976403
$ws.Range("H6").Value = 1.233370515372144
$ws.Range("H7").Value = 1.231176922344901
$ws.Range("H8").Value = 1.222617728650292
$ws.Range("H9").Value = 1.209860144005603
$ws.Range("H10").Value = 1.202942357298411
$ws.Range("H11").Value = 1.200328552392619
$ws.Range("H12").Value = 1.199415443178538
$ws.Range("H13").Value = 1.199608686824334
$ws.Range("H14").Value = 1.200251893416507
$ws.Range("H15").Value = 1.200655862853424
$ws.Range("H16").Value = 1.203123905093435
$ws.Range("H17").Value = 1.204774529330308
$ws.Range("H18").Value = 1.205774104339724
$ws.Range("H19").Value = 1.206121161126802
$ws.Range("H20").Value = 1.204593624029059
$ws.Range("H21").Value = 1.200060886768199
$ws.Range("H22").Value = 1.197545451751296
$ws.Range("H23").Value = 1.198847082452417
$ws.Range("H24").Value = 1.204675253666565
$ws.Range("H25").Value = 1.212880272640845
$ws.Range("K2").Value = 0.4897642838174079
$ws.Range("K3").Value = 0.4533821092130381
$ws.Range("K4").Value = 0.4312525970329659
$ws.Range("K5").Value = 0.4222874953840687
$ws.Range("K6").Value = 0.4208020460889941
$ws.Range("K7").Value = 0.4311314762305187
$ws.Range("K8").Value = 0.4771763952184926
$ws.Range("K9").Value = 0.5691260797274253
$ws.Range("K10").Value = 0.6376918255818111
$ws.Range("K11").Value = 0.6691042216453411
$ws.Range("K12").Value = 0.6810310170907314
$ws.Range("K13").Value = 0.6784609666860888
$ws.Range("K14").Value = 0.6700848143333644
$ws.Range("K15").Value = 0.6649582874396742
$ws.Range("K16").Value = 0.6356433874494485
$ws.Range("K17").Value = 0.6177161991874982
$ws.Range("K18").Value = 0.6074258175073624
$ws.Range("K19").Value = 0.6039452618657606
$ws.Range("K20").Value = 0.6196224215545101
$ws.Range("K21").Value = 0.6725442375998227
$ws.Range("K22").Value = 0.7073157946697108
$ws.Range("K23").Value = 0.6887407975471262
$ws.Range("K24").Value = 0.6187605679514832
$ws.Range("K25").Value = 0.5440739039493678
$ws.Range("L2").Value = 0.2076874156526998
$ws.Range("L3").Value = 0.2008295540185259
$ws.Range("L4").Value = 0.196734486755318
$ws.Range("L5").Value = 0.1950948371984111
$ws.Range("L6").Value = 0.1948243346214582
$ws.Range("L7").Value = 0.1967122559030656
$ws.Range("L8").Value = 0.2052988311048694
$ws.Range("L9").Value = 0.2230546989700315
$ws.Range("L10").Value = 0.2366606187650149
$ws.Range("L11").Value = 0.2429724328912783
$ws.Range("L12").Value = 0.2453801499539594
$ws.Range("L13").Value = 0.244860824094232
$ws.Range("L14").Value = 0.2431701652805742
$ws.Range("L15").Value = 0.2421368751062971
$ws.Range("L16").Value = 0.236250587324335
$ws.Range("L17").Value = 0.2326708699620781
$ws.Range("L18").Value = 0.230623435262828
$ws.Range("L19").Value = 0.2299321893193991
$ws.Range("L20").Value = 0.2330507445297627
$ws.Range("L21").Value = 0.2436662761718651
$ws.Range("L22").Value = 0.2507065362211875
$ws.Range("L23").Value = 0.2469396584669568
$ws.Range("L24").Value = 0.2328789702422114
$ws.Range("L25").Value = 0.2181528971193245
$ws.Range("N2").Value = 2.206381522087909
$ws.Range("N3").Value = 2.22848029150143
$ws.Range("N4").Value = 2.242758706423654
$ws.Range("N5").Value = 2.248755669054308
$ws.Range("N6").Value = 2.249762232827869
$ws.Range("N7").Value = 2.242838861341276
$ws.Range("N8").Value = 2.213853847843321
$ws.Range("N9").Value = 2.162647277537499
$ws.Range("N10").Value = 2.128461268563651
$ws.Range("N11").Value = 2.113655475020426
$ws.Range("N12").Value = 2.108156201289521
$ws.Range("N13").Value = 2.109335791136928
$ws.Range("N14").Value = 2.113200893992374
$ws.Range("N15").Value = 2.115582366722677
$ws.Range("N16").Value = 2.12944387742019
$ws.Range("N17").Value = 2.138138530554734
$ws.Range("N18").Value = 2.143209635010134
$ws.Range("N19").Value = 2.144938671296259
$ws.Range("N20").Value = 2.137205706525798
$ws.Range("N21").Value = 2.112062704738165
$ws.Range("N22").Value = 2.096256139040612
$ws.Range("N23").Value = 2.104635087891051
$ws.Range("N24").Value = 2.1376272106457
$ws.Range("N25").Value = 2.175896829143909
